# Weekly refresh of the "Fruta, Terminal La Palmera de La Serena - Membrillo" sheet.
# Rows 2-8 and 11-12 are rewritten with the latest weekly market data; rows 9-10
# (the most recent week already on file) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for each affected row, keyed by destination row number.
# Columns: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
#          O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44309; L = "Especial"; M = 20; N = 305000; O = 310000; P = 307500; R = "Provincia de Cachapoal"; S = 683 },
    @{ Row = 3;  D = 44309; L = "Primera";  M = 20; N = 285000; O = 290000; P = 287500; R = "Provincia de Cachapoal"; S = 639 },
    @{ Row = 4;  D = 44309; L = "Segunda";  M = 20; N = 255000; O = 260000; P = 257500; R = "Provincia de Cachapoal"; S = 572 },
    @{ Row = 5;  D = 44273; L = "Especial"; M = 10; N = 255000; O = 260000; P = 257500; R = "Región de O'Higgins";   S = 572 },
    @{ Row = 6;  D = 44273; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región de O'Higgins";   S = 506 },
    @{ Row = 7;  D = 44295; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región Metropolitana";  S = 506 },
    @{ Row = 8;  D = 44295; L = "Segunda";  M = 16; N = 195000; O = 200000; P = 197500; R = "Región Metropolitana";  S = 439 },
    @{ Row = 11; D = 44294; L = "Primera";  M = 20; N = 225000; O = 230000; P = 227500; R = "Región Metropolitana";  S = 506 },
    @{ Row = 12; D = 44294; L = "Segunda";  M = 16; N = 195000; O = 200000; P = 197500; R = "Región Metropolitana";  S = 439 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
}
